$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = "NIG(1.092629183047433, 0.6757123593054357, 1.2104671275467238, 3.7924345806441857)"
$ws.Range("C2").Value = "NIG(1.7288981085724426, 1.2400512705037015, 3.600629997067743, 6.303791028170348)"
$ws.Range("D2").Value = "JSU(-1.0834270695418797, 1.1147815040763032, -0.09648624342902937, 2.6439773442607537)"
$ws.Range("E2").Value = "JSU(-0.9462123012688848, 1.137637810872195, 3.48742340195479, 5.201982745315803)"
